$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("O2").Value = "2022-09-13 21:00:31"

# Row 3
$ws.Range("O3").Value = "2022-09-13 21:00:31"

# Row 4
$ws.Range("O4").Value = "2022-09-13 21:00:31"

# Row 5
$ws.Range("O5").Value = "2022-09-13 21:00:31"

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "6568452"
$ws.Range("B6").Value = "Super Soft Premium Mandel feucht 4x  50ST"
$ws.Range("C6").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/super-soft-premium-mandel-feucht/p/6568452"
$ws.Range("D6").Value = "4x 50ST"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 3.5
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "7.65"
$ws.Range("I6").Value = "0.04/1ST"
$ws.Range("J6").Value = "Preis pro 1 Stück"
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "0.04"
$ws.Range("L6").Value = "1ST"
$ws.Range("M6").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Range("N6").Value = "Super Soft Premium Mandel feucht 4x  50ST 35% Aktion 7.65 Schweizer Franken statt 11.80 Schweizer Franken"
$ws.Range("O6").Value = "2022-09-13 21:00:31"

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "6497242"
$ws.Range("B7").Value = "Super Soft WC-Papier Sensation, 4 Rollen 3-lagig"
$ws.Range("C7").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/super-soft-wc-papier-sensation-4-rollen-3-lagig/p/6497242"
$ws.Range("D7").Value = "4Rol"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 4.5
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "3.25"
$ws.Range("I7").Value = "0.81/1Rol"
$ws.Range("J7").Value = "Preis pro 1 Rolle"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "0.81"
$ws.Range("L7").Value = "1Rol"
$ws.Range("M7").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N7").Value = "Super Soft WC-Papier Sensation, 4 Rollen 3-lagig 3.25 Schweizer Franken"
$ws.Range("O7").Value = "2022-09-13 21:00:31"

# Row 8
$ws.Range("O8").Value = "2022-09-13 21:00:31"

# Row 9
$ws.Range("O9").Value = "2022-09-13 21:00:31"

# Row 10
$ws.Range("O10").Value = "2022-09-13 21:00:31"

# Row 11
$ws.Range("O11").Value = "2022-09-13 21:00:31"

# Row 12
$ws.Range("O12").Value = "2022-09-13 21:00:31"

# Row 13
$ws.Range("O13").Value = "2022-09-13 21:00:31"

# Row 14
$ws.Range("O14").Value = "2022-09-13 21:00:31"

# Row 15
$ws.Range("O15").Value = "2022-09-13 21:00:31"

# Row 16
$ws.Range("O16").Value = "2022-09-13 21:00:31"

# Row 17
$ws.Range("O17").Value = "2022-09-13 21:00:31"

# Row 18
$ws.Range("O18").Value = "2022-09-13 21:00:31"

# Row 19
$ws.Range("O19").Value = "2022-09-13 21:00:31"

# Row 20
$ws.Range("O20").Value = "2022-09-13 21:00:31"

# Row 21
$ws.Range("O21").Value = "2022-09-13 21:00:31"

# Row 22
$ws.Range("O22").Value = "2022-09-13 21:00:31"

# Row 23
$ws.Range("O23").Value = "2022-09-13 21:00:31"

# Row 24
$ws.Range("O24").Value = "2022-09-13 21:00:31"

# Row 25
$ws.Range("O25").Value = "2022-09-13 21:00:31"

# Row 26
$ws.Range("O26").Value = "2022-09-13 21:00:31"

# Row 27
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = "3874909"
$ws.Range("B27").Value = "Oecoplan Papiertaschentücher Special-Edition Calendula 30x10 Stück"
$ws.Range("C27").Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/oecoplan-papiertaschentuecher-special-edition-calendula-30x10-stueck/p/3874909"
$ws.Range("D27").Value = "30ST"
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = "Coop"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "3.65"
$ws.Range("I27").Value = "0.12/1ST"
$ws.Range("J27").Value = "Preis pro 1 Stück"
$ws.Range("K27").NumberFormat = "@"
$ws.Range("K27").Value = "0.12"
$ws.Range("L27").Value = "1ST"
$ws.Range("M27").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Range("N27").Value = "Oecoplan Papiertaschentücher Special-Edition Calendula 30x10 Stück 20% Aktion 3.65 Schweizer Franken statt 4.60 Schweizer Franken"
$ws.Range("O27").Value = "2022-09-13 21:00:31"

# Row 28
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = "6996030"
$ws.Range("B28").Value = "Tela Viva Haushaltspapier 3-lagig 4 Rollen"
$ws.Range("C28").Value = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/tela-viva-haushaltspapier-3-lagig-4-rollen/p/6996030"
$ws.Range("D28").Value = "200BLT"
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = "Tela"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "5.95"
$ws.Range("I28").Value = ""
$ws.Range("J28").Value = ""
$ws.Range("K28").Value = ""
$ws.Range("L28").Value = ""
$ws.Range("M28").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
$ws.Range("N28").Value = "Tela Viva Haushaltspapier 3-lagig 4 Rollen 5.95 Schweizer Franken"
$ws.Range("O28").Value = "2022-09-13 21:00:31"

# Row 29
$ws.Range("O29").Value = "2022-09-13 21:00:31"

# Row 30
$ws.Range("O30").Value = "2022-09-13 21:00:31"

# Row 31
$ws.Range("O31").Value = "2022-09-13 21:00:31"

# Row 32
$ws.Range("O32").Value = "2022-09-13 21:00:31"

# Row 33
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "6283677"
$ws.Range("B33").Value = "Oecoplan Goldmelisse blau 3-lagig 32 Rollen"
$ws.Range("C33").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/oecoplan-goldmelisse-blau-3-lagig-32-rollen/p/6283677"
$ws.Range("D33").Value = "32Rol"
$ws.Range("E33").Value = 1
$ws.Range("F33").Value = 5
$ws.Range("G33").Value = "Coop"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "14.80"
$ws.Range("I33").Value = "0.46/1Rol"
$ws.Range("J33").Value = "Preis pro 1 Rolle"
$ws.Range("K33").NumberFormat = "@"
$ws.Range("K33").Value = "0.46"
$ws.Range("L33").Value = "1Rol"
$ws.Range("M33").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N33").Value = "Oecoplan Goldmelisse blau 3-lagig 32 Rollen 30% Aktion 14.80 Schweizer Franken statt 21.20 Schweizer Franken"
$ws.Range("O33").Value = "2022-09-13 21:00:31"

# Row 34
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = "6636712"
$ws.Range("B34").Value = "Pampers Coconut Pure 42 Feuchttücher"
$ws.Range("C34").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/pampers-coconut-pure-42-feuchttuecher/p/6636712"
$ws.Range("D34").Value = "42ST"
$ws.Range("E34").Value = 2
$ws.Range("F34").Value = 3.5
$ws.Range("G34").Value = "Pampers"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "4.95"
$ws.Range("I34").Value = "0.12/1ST"
$ws.Range("K34").NumberFormat = "@"
$ws.Range("K34").Value = "0.12"
$ws.Range("N34").Value = "Pampers Coconut Pure 42 Feuchttücher 4.95 Schweizer Franken"
$ws.Range("O34").Value = "2022-09-13 21:00:31"

# Row 35
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = "6727164"
$ws.Range("B35").Value = "Wetties Allzwecktücher 80Stück"
$ws.Range("C35").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/wetties-allzwecktuecher-80stueck/p/6727164"
$ws.Range("D35").Value = "80ST"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "3.95"
$ws.Range("I35").Value = "0.05/1ST"
$ws.Range("J35").Value = "Preis pro 1 Stück"
$ws.Range("K35").NumberFormat = "@"
$ws.Range("K35").Value = "0.05"
$ws.Range("L35").Value = "1ST"
$ws.Range("M35").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Range("N35").Value = "Wetties Allzwecktücher 80Stück 3.95 Schweizer Franken"
$ws.Range("O35").Value = "2022-09-13 21:00:31"

